$d = $word.ActiveDocument

function Get-ParaIndexByText($text) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $t = $d.Paragraphs($i).Range.Text
        $t = $t.TrimEnd([char]13, [char]7)
        if ($t -eq $text) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# 1. "state of the art" -> "state-of-the-art" in the background/problem
#    paragraph (keeps the existing run formatting, sz=24 / szCs=24)
# ---------------------------------------------------------------------------
$bgFull = "The development of sophisticated imaging technologies and machine learning presents a possible alternative to the poor performance of existing methods for the identification of prostate cancer. Utilising these technologies can lead to better patient outcomes by enhancing diagnostic accuracy and personalising treatment approaches. To overcome current obstacles and fully realise the potential of these state of the art instruments in the treatment of prostate cancer, more research and innovation in this area are imperative."
$bgIdx = Get-ParaIndexByText($bgFull)
if ($bgIdx -gt 0) {
    $p = $d.Paragraphs($bgIdx)
    $pText = $p.Range.Text
    $marker = "personalising treatment"
    $splitPoint = $pText.IndexOf($marker) + $marker.Length
    $pStart = $p.Range.Start
    $run2 = $d.Range($pStart + $splitPoint, $p.Range.End - 1)
    $run2.Text = " approaches. To overcome current obstacles and fully realise the potential of these state-of-the-art instruments in the treatment of prostate cancer, more research and innovation in this area are imperative."
}

# ---------------------------------------------------------------------------
# 2. "Specific problem my project is addressing" -> real problem-statement
#    text, followed by three brand new paragraphs of body text.
# ---------------------------------------------------------------------------
$idx = Get-ParaIndexByText("Specific problem my project is addressing")
if ($idx -gt 0) {
    $p = $d.Paragraphs($idx)
    $p.Range.Text = "Prostate cancer is a major worldwide health concern, as it is one of the most often diagnosed cancers among the men and the leading cause of cancer related death. For successful therapy and better patient outcomes, clinically significant prostate cancer lesions must be identified early and accurate. Due to the poor sensitivity and specificity of traditional diagnostic techniques such as digital rectal exams (DRE) and prostate specific antigen (PSA) testing, there is a risk of overdiagnosis, overtreatment, or missing diagnoses."

    $p.Range.InsertParagraphAfter()
    $idx = $idx + 1
    $p = $d.Paragraphs($idx)
    $p.Range.Text = "Multiparametric magnetic resonance imaging (mpMRI) improves the capacity to detect and characterise prostate cancer by providing both anatomical and functional imaging. The assessment of prostate lesions clinical significance (ClinSig) using mpMRI picture interpretation is still difficult and heavily reliant on radiologists’ skill, which often results in inter observer variability."

    $p.Range.InsertParagraphAfter()
    $idx = $idx + 1
    $p = $d.Paragraphs($idx)
    $p.Range.Text = "Machine learning offers a promising answer to these problems by automating the interpretation of mpMRI images and predicting the ClinSig score of prostate lesions. Large amounts of imaging data can be processed by ML models, especially deep learning approaches, which can then be used to spot subtle patterns that human observers might miss, enhancing diagnostic consistency and accuracy."

    $p.Range.InsertParagraphAfter()
    $idx = $idx + 1
    $p = $d.Paragraphs($idx)
    $p.Range.Text = "This project focuses on creating and deploying a machine learning model for predicting the ClinSig score of prostate lesions based on T2 weighted mpMRI images."
}

# ---------------------------------------------------------------------------
# 3. "why this study is important and its potential impact" -> justification
#    text
# ---------------------------------------------------------------------------
$idx = Get-ParaIndexByText("why this study is important and its potential impact")
if ($idx -gt 0) {
    $p = $d.Paragraphs($idx)
    $p.Range.Text = "This study is primarily justified by the possiblility that it may greatly enhance the diagnosis of prostate cancer by offering a more precise and reliable way to predict the ClinSig score of prostate lesions. By lowering diagnostic mistakes and inter observer variability, this automated method can assist radiologists in making more informed treatment decisions. Furthermore, the effective diagnostic process can be improved by integrating machine learning models into clinical processes. This guarantees prompt and suitable interventions, which are critical for improving patient outcomes. This project aims to lessen the burden of prostate cancer on healthcare systems and contribute to personalised treatment regimens by increasing regimens by increasing diagnosis accuracy and consistency."
}

# ---------------------------------------------------------------------------
# 4. "Research questions that I want to answer by this project" -> research
#    question, now a numbered list item (numId 7 in the target document).
# ---------------------------------------------------------------------------
$idx = Get-ParaIndexByText("Research questions that I want to answer by this project")
if ($idx -gt 0) {
    $p = $d.Paragraphs($idx)
    $p.Range.Text = "How accurately can a convolutional neural network (CNN) model predict the ClinSig score of prostate lesions form the T2 weighted mpMRI images?"
    $p.Style = "ListParagraph"
    $p.Range.ListFormat.ApplyNumberDefault()
}

# ---------------------------------------------------------------------------
# 5. "Project aims and objectives" -> aims paragraph, followed by four new
#    numbered objective paragraphs (numId 8 in the target document).
# ---------------------------------------------------------------------------
$idx = Get-ParaIndexByText("Project aims and objectives")
if ($idx -gt 0) {
    $p = $d.Paragraphs($idx)
    $p.Range.Text = "This project’s main goal is to create and verify a machine learning model that may be used to reliably predict prostate lesion clinical significance (ClinSig) scores using multiparametric magnetic resonance images (mpMRI) data. This goal will be met through the following specific objectives."

    # First list item gets a brand-new numbered list (ApplyNumberDefault);
    # the following three simply continue that same list because they are
    # created via InsertParagraphAfter on an already-numbered paragraph, so
    # they inherit its pPr (pStyle + numPr) including the numId.
    $p.Range.InsertParagraphAfter()
    $idx = $idx + 1
    $p = $d.Paragraphs($idx)
    $p.Range.Text = "Investigate the machine learning methods that are currently being utilised to categorise and predict prostate cancer."
    $p.Style = "ListParagraph"
    $p.Range.ListFormat.ApplyNumberDefault()

    $p.Range.InsertParagraphAfter()
    $idx = $idx + 1
    $p = $d.Paragraphs($idx)
    $p.Range.Text = "Load and prepare the Prostatex dataset images (T2 weighted images) specified in the detailed description from the Prostatex challenge."

    $p.Range.InsertParagraphAfter()
    $idx = $idx + 1
    $p = $d.Paragraphs($idx)
    $p.Range.Text = "Using the prepared data, create a convolutional neural network model to predict the ClinSig score of prostate lesions"

    $p.Range.InsertParagraphAfter()
    $idx = $idx + 1
    $p = $d.Paragraphs($idx)
    $p.Range.Text = "Evaluate the created CNN models performance in comparison to current classifier models and conventional diagnostic techniques"
}

# ---------------------------------------------------------------------------
# 6. Drop the stale w:lastRenderedPageBreak markers on the two headings that
#    no longer start a fresh page ("2.4 Related work" and "4.3 Comparison
#    with the existing methods"). Re-asserting the text rebuilds the run and
#    drops the obsolete rendering bookmark.
# ---------------------------------------------------------------------------
foreach ($headingText in @("2.4 Related work", "4.3 Comparison with the existing methods")) {
    $hIdx = Get-ParaIndexByText($headingText)
    if ($hIdx -gt 0) {
        $hp = $d.Paragraphs($hIdx)
        $hp.Range.Text = $headingText
    }
}
